# Update parameters output with new equation (adds beta1/beta2 rows,
# shifts gamma1/gamma2/nu/rho/phi1 down, updates their values, and
# appends a new phi2 row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: beta1
$ws.Range("B2").Value = "beta1"
$ws.Range("C2").Value = 0.97
$ws.Range("D2").Value = 0.97
$ws.Range("E2").Value = 0.84
$ws.Range("F2").Value = 1.06

# Row 3: beta2
$ws.Range("B3").Value = "beta2"
$ws.Range("C3").Value = 0.51
$ws.Range("D3").Value = 0.58
$ws.Range("E3").Value = 0.04
$ws.Range("F3").Value = 1.43

# Row 4: gamma1
$ws.Range("A4").Value = 17
$ws.Range("B4").Value = "gamma1"
$ws.Range("C4").Value = 0.25
$ws.Range("D4").Value = 0.25
$ws.Range("E4").Value = 0.15
$ws.Range("F4").Value = 0.34

# Row 5: gamma2
$ws.Range("A5").Value = 17
$ws.Range("B5").Value = "gamma2"
$ws.Range("C5").Value = 0.25
$ws.Range("D5").Value = 0.25
$ws.Range("E5").Value = 0.15
$ws.Range("F5").Value = 0.35

# Row 6: nu
$ws.Range("A6").Value = 17
$ws.Range("B6").Value = "nu"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0

# Row 7: rho
$ws.Range("A7").Value = 17
$ws.Range("B7").Value = "rho"
$ws.Range("C7").Value = 0.33
$ws.Range("D7").Value = 0.34
$ws.Range("E7").Value = 0.06
$ws.Range("F7").Value = 0.71

# Row 8: phi1 (new row, shifted from old row 6)
$ws.Range("A8").Value = 17
$ws.Range("B8").Value = "phi1"
$ws.Range("C8").Value = 0.49
$ws.Range("D8").Value = 0.49
$ws.Range("E8").Value = 0.34
$ws.Range("F8").Value = 0.7

# Row 9: phi2 (brand new row)
$ws.Range("A9").Value = 17
$ws.Range("B9").Value = "phi2"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
